$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 3.156260333333333
$ws.Range("N2").Value = 9.468781
$ws.Range("O2").Value = 0.3177865155521988
$ws.Range("P2").Value = 0.3177865155521988
$ws.Range("Q2").Value = 1.432548710878444
$ws.Range("R2").Value = 12.892938397906
$ws.Range("S2").Value = 0.3177865155521988
$ws.Range("T2").Value = 0.3177865155521988

# Row 3
$ws.Range("M3").Value = 2.821123666666667
$ws.Range("N3").Value = 8.463371
$ws.Range("O3").Value = 0.2840434454990065
$ws.Range("P3").Value = 0.2840434454990065
$ws.Range("S3").Value = 0.2840434454990065
$ws.Range("T3").Value = 0.2840434454990065

# Row 4
$ws.Range("M4").Value = 3.954630666666667
$ws.Range("N4").Value = 11.863892
$ws.Range("O4").Value = 0.3981700389487947
$ws.Range("P4").Value = 0.3981700389487947
$ws.Range("Q4").Value = 1.794909312043556
$ws.Range("R4").Value = 16.154183808392
$ws.Range("S4").Value = 0.3981700389487947
$ws.Range("T4").Value = 0.3981700389487947
